$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base timer frequency formula (200 -> 240 MHz); dependents recalc automatically
$ws.Range("B4").Formula = "=240*1000*1000"

# Update the active cell selection on the sheet view
$ws.Range("I10").Select()
